$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift existing data rows 2-21 down by 7 (become rows 9-28). ---
# Walk bottom-up so a row is always fully written with its original values
# before the slot above it (which holds the same data) gets overwritten.
$ws.Range("A28").Value = 0.0099265603348612
$ws.Range("B28").Value = -0.0109955742955207
$ws.Range("C28").Value = -0.0491746515035629
$ws.Range("A27").Value = -0.0064140851609408
$ws.Range("B27").Value = 0.0251981914043426
$ws.Range("C27").Value = 0.0591012127697467
$ws.Range("A26").Value = -0.0339030213654041
$ws.Range("B26").Value = -0.0348193198442459
$ws.Range("C26").Value = 0.0242818929255008
$ws.Range("A25").Value = 0.0025961773935705
$ws.Range("B25").Value = -0.0377209298312664
$ws.Range("C25").Value = -0.1343903541564941
$ws.Range("A24").Value = -0.012980886735022
$ws.Range("B24").Value = 0.1769981980323791
$ws.Range("C24").Value = 0.0358883328735828
$ws.Range("A23").Value = 0.0079412478953599
$ws.Range("B23").Value = -0.0103847095742821
$ws.Range("C23").Value = 0.0820086598396301
$ws.Range("A22").Value = -0.0311541277915239
$ws.Range("B22").Value = 0.040775254368782
$ws.Range("C22").Value = -0.4699080884456634
$ws.Range("A21").Value = 0.001527163083665
$ws.Range("B21").Value = 0.3072652220726013
$ws.Range("C21").Value = 0.0675006061792373
$ws.Range("A20").Value = -0.1476766765117645
$ws.Range("B20").Value = 0.0491746515035629
$ws.Range("C20").Value = 0.3875939846038818
$ws.Range("A19").Value = -0.1258382350206375
$ws.Range("B19").Value = -0.4847215712070465
$ws.Range("C19").Value = -0.2483167201280594
$ws.Range("A18").Value = -0.2658790946006775
$ws.Range("B18").Value = -0.8356636762619019
$ws.Range("C18").Value = -1.123686671257019
$ws.Range("A17").Value = -0.3932445049285888
$ws.Range("B17").Value = -0.6039929986000061
$ws.Range("C17").Value = 1.466076612472534
$ws.Range("A16").Value = 0.3782783150672912
$ws.Range("B16").Value = -0.1327104717493057
$ws.Range("C16").Value = 0.1701259762048721
$ws.Range("A15").Value = 1.769065737724304
$ws.Range("B15").Value = -1.171639561653137
$ws.Range("C15").Value = -0.7515169382095337
$ws.Range("A14").Value = -1.828319668769836
$ws.Range("B14").Value = 1.374752283096314
$ws.Range("C14").Value = 2.351678371429444
$ws.Range("A13").Value = 0.1690569519996643
$ws.Range("B13").Value = -0.0616973899304866
$ws.Range("C13").Value = -0.9847147464752196
$ws.Range("A12").Value = -0.2105957865715026
$ws.Range("B12").Value = 0.9750936627388
$ws.Range("C12").Value = -1.285871386528015
$ws.Range("A11").Value = -0.3101668357849121
$ws.Range("B11").Value = -0.0274889357388019
$ws.Range("C11").Value = -0.0977384373545646
$ws.Range("A10").Value = -0.011148290708661
$ws.Range("B10").Value = -0.4109596014022827
$ws.Range("C10").Value = 0.3446807265281677
$ws.Range("A9").Value = -0.294895201921463
$ws.Range("B9").Value = 0.9447031021118164
$ws.Range("C9").Value = 0.0591012127697467

# --- Step 2: write the 7 newly-inserted rows at the top (rows 2-8). ---
$ws.Range("A2").Value = -0.005192354787141
$ws.Range("B2").Value = 0.0755945742130279
$ws.Range("C2").Value = 0.0082466807216405
$ws.Range("A3").Value = 0.00167987938039
$ws.Range("B3").Value = 0.0209221355617046
$ws.Range("C3").Value = -0.0082466807216405
$ws.Range("A4").Value = 0.0074830991216003
$ws.Range("B4").Value = -0.0106901414692401
$ws.Range("C4").Value = 0.0113010071218013
$ws.Range("A5").Value = -0.0024434609804302
$ws.Range("B5").Value = 0.0310014113783836
$ws.Range("C5").Value = 0.0187841057777404
$ws.Range("A6").Value = -0.0039706239476799
$ws.Range("B6").Value = 0.0114537235349416
$ws.Range("C6").Value = -0.0603229440748691
$ws.Range("A7").Value = -0.0891863256692886
$ws.Range("B7").Value = 0.2981022298336029
$ws.Range("C7").Value = 0.0048869219608604
$ws.Range("A8").Value = -0.2064724564552307
$ws.Range("B8").Value = 0.7906123399734497
$ws.Range("C8").Value = 0.058184914290905

# --- Step 3: append the 3 new rows at the bottom (rows 29-31). ---
$ws.Range("A29").Value = -0.0039706239476799
$ws.Range("B29").Value = 0.0164933614432811
$ws.Range("C29").Value = -0.0076358155347406
$ws.Range("A30").Value = -0.0065668015740811
$ws.Range("B30").Value = -0.0163406450301408
$ws.Range("C30").Value = 0.0088575463742017
$ws.Range("A31").Value = -0.0093156946823
$ws.Range("B31").Value = -0.0369573459029197
$ws.Range("C31").Value = -0.0157297793775796

Write-Output "done"
